$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nts"
$ws.Range("C2").Value = "Sort1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.729615
$ws.Range("H2").Value = 8.188845000000001
$ws.Range("I2").Value = 0.8986540674935011
$ws.Range("J2").Value = 0.9300737996692831
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.9703876666666668
$ws.Range("N2").Value = 2.911163
$ws.Range("O2").Value = 0.03945299285965207
$ws.Range("P2").Value = 0.04754668824173519
$ws.Range("Q2").Value = 2.648784730748334
$ws.Range("R2").Value = 23.839062576735
$ws.Range("S2").Value = 0.03545459250811839
$ws.Range("T2").Value = 0.04422192899468147

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nts"
$ws.Range("C3").Value = "Sort1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.729615
$ws.Range("H3").Value = 8.188845000000001
$ws.Range("I3").Value = 0.8986540674935011
$ws.Range("J3").Value = 0.9300737996692831
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.041192666666667
$ws.Range("N3").Value = 3.123578
$ws.Range("O3").Value = 0.04233170747586662
$ws.Range("P3").Value = 0.05101596487889641
$ws.Range("Q3").Value = 2.842055120823334
$ws.Range("R3").Value = 25.57849608741
$ws.Range("S3").Value = 0.03804156110713259
$ws.Range("T3").Value = 0.04744861229870988

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Nts"
$ws.Range("C4").Value = "Sort1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.729615
$ws.Range("H4").Value = 8.188845000000001
$ws.Range("I4").Value = 0.8986540674935011
$ws.Range("J4").Value = 0.9300737996692831
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.170211333333333
$ws.Range("N4").Value = 15.510634
$ws.Range("O4").Value = 0.2102049704707969
$ws.Range("P4").Value = 0.2533280614069559
$ws.Range("Q4").Value = 14.11268640863667
$ws.Range("R4").Value = 127.01417767773
$ws.Range("S4").Value = 0.188901551720933
$ws.Range("T4").Value = 0.2356137926356209

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Nts"
$ws.Range("C5").Value = "Sort1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.729615
$ws.Range("H5").Value = 8.188845000000001
$ws.Range("I5").Value = 0.8986540674935011
$ws.Range("J5").Value = 0.9300737996692831
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.853575333333334
$ws.Range("N5").Value = 14.560726
$ws.Range("O5").Value = 0.1973315197085667
$ws.Range("P5").Value = 0.2378136503161547
$ws.Range("Q5").Value = 13.24839203349667
$ws.Range("R5").Value = 119.23552830147
$ws.Range("S5").Value = 0.1773327728307775
$ws.Range("T5").Value = 0.2211842453627682

$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Nts"
$ws.Range("C6").Value = "Sort1"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.729615
$ws.Range("H6").Value = 8.188845000000001
$ws.Range("I6").Value = 0.8986540674935011
$ws.Range("J6").Value = 0.9300737996692831
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 12.56068
$ws.Range("N6").Value = 25.12136
$ws.Range("O6").Value = 0.5106788094851177
$ws.Range("P6").Value = 0.4102956351562577
$ws.Range("Q6").Value = 34.28582053820001
$ws.Range("R6").Value = 205.7149232292
$ws.Range("S6").Value = 0.4589235893265397
$ws.Range("T6").Value = 0.3816052203775025

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Nts"
$ws.Range("C7").Value = "Sort1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.307833
$ws.Range("H7").Value = 0.615666
$ws.Range("I7").Value = 0.1013459325064989
$ws.Range("J7").Value = 0.06992620033071682
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.9703876666666668
$ws.Range("N7").Value = 2.911163
$ws.Range("O7").Value = 0.03945299285965207
$ws.Range("P7").Value = 0.04754668824173519
$ws.Range("Q7").Value = 0.298717346593
$ws.Range("R7").Value = 1.792304079558
$ws.Range("S7").Value = 0.003998400351533681
$ws.Range("T7").Value = 0.003324759247053713

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Nts"
$ws.Range("C8").Value = "Sort1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.307833
$ws.Range("H8").Value = 0.615666
$ws.Range("I8").Value = 0.1013459325064989
$ws.Range("J8").Value = 0.06992620033071682
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.041192666666667
$ws.Range("N8").Value = 3.123578
$ws.Range("O8").Value = 0.04233170747586662
$ws.Range("P8").Value = 0.05101596487889641
$ws.Range("Q8").Value = 0.320513462158
$ws.Range("R8").Value = 1.923080772948
$ws.Range("S8").Value = 0.004290146368734032
$ws.Range("T8").Value = 0.003567352580186524

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Nts"
$ws.Range("C9").Value = "Sort1"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.307833
$ws.Range("H9").Value = 0.615666
$ws.Range("I9").Value = 0.1013459325064989
$ws.Range("J9").Value = 0.06992620033071682
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.170211333333333
$ws.Range("N9").Value = 15.510634
$ws.Range("O9").Value = 0.2102049704707969
$ws.Range("P9").Value = 0.2533280614069559
$ws.Range("Q9").Value = 1.591561665374
$ws.Range("R9").Value = 9.549369992244001
$ws.Range("S9").Value = 0.02130341874986397
$ws.Range("T9").Value = 0.01771426877133493

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Nts"
$ws.Range("C10").Value = "Sort1"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.307833
$ws.Range("H10").Value = 0.615666
$ws.Range("I10").Value = 0.1013459325064989
$ws.Range("J10").Value = 0.06992620033071682
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.853575333333334
$ws.Range("N10").Value = 14.560726
$ws.Range("O10").Value = 0.1973315197085667
$ws.Range("P10").Value = 0.2378136503161547
$ws.Range("Q10").Value = 1.494090655586
$ws.Range("R10").Value = 8.964543933516001
$ws.Range("S10").Value = 0.01999874687778926
$ws.Range("T10").Value = 0.01662940495338647

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Nts"
$ws.Range("C11").Value = "Sort1"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.307833
$ws.Range("H11").Value = 0.615666
$ws.Range("I11").Value = 0.1013459325064989
$ws.Range("J11").Value = 0.06992620033071682
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 12.56068
$ws.Range("N11").Value = 25.12136
$ws.Range("O11").Value = 0.5106788094851177
$ws.Range("P11").Value = 0.4102956351562577
$ws.Range("Q11").Value = 3.86659180644
$ws.Range("R11").Value = 15.46636722576
$ws.Range("S11").Value = 0.05175522015857793
$ws.Range("T11").Value = 0.02869041477875518

